$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Price cells are stored as text (e.g. "65.548.18", "1.00") in the source sheet,
# so we force a text NumberFormat while writing the value, then restore the
# cell style to Normal so no stray formatting is left behind on save.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "65.548.18"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +1.88%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.645.17"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "604.55"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "156.25"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +2.46%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.53%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.643.40"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.49%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.123"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +7.38%  "
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("E13").Value = "  +1.54%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "29.68"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +5.70%  "
$ws.Range("E15").Value = "  +13.77%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.123.54"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.60%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "65.285.64"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.69%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.642.83"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.84%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "12.61"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.38%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "4.87"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.94%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "357.66"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +1.95%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "7.43"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +4.62%  "
$ws.Range("E23").Value = "  -0.01%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "69.52"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +2.64%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "1.70"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.38"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("E27").Value = "  +15.36%  "
$ws.Range("E28").Value = "  -3.11%  "
$ws.Range("E29").Value = "  +1.80%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "8.10"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +4.38%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "526.78"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -4.97%  "
$ws.Range("E34").Value = "  -3.38%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "5.52"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.37%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "6.33"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.96%  "
$ws.Range("E37").Value = "  +1.86%  "
$ws.Range("E38").Value = "  +2.80%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "161.47"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.79%  "
$ws.Range("E40").Value = "  -1.38%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  -0.01%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "41.91"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +3.90%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "165.04"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("E45").Value = "  +0.51%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.35"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +6.05%  "
$ws.Range("E47").Value = "  +2.96%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "22.84"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("E50").Value = "  +3.00%  "
$ws.Range("E51").Value = "  +0.42%  "
